$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range keeps its original text formatting so that
# numeric-looking values (e.g. "1.00", "0.997") are not auto-converted
# into Excel numbers by the COM layer and lose formatting.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '54.162.33'
$ws.Range('E2').Value = '  -10.71%  '
$ws.Range('D3').Value = '2.352.02'
$ws.Range('E3').Value = '  -19.12%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '446.49'
$ws.Range('E5').Value = '  -15.64%  '
$ws.Range('D6').Value = '128.31'
$ws.Range('E6').Value = '  -11.03%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.477'
$ws.Range('E8').Value = '  -14.32%  '
$ws.Range('D9').Value = '2.360.32'
$ws.Range('E9').Value = '  -19.07%  '
$ws.Range('D10').Value = '5.40'
$ws.Range('E10').Value = '  -10.69%  '
$ws.Range('D11').Value = '0.0923'
$ws.Range('E11').Value = '  -15.14%  '
$ws.Range('D12').Value = '0.310'
$ws.Range('E12').Value = '  -14.54%  '
$ws.Range('E13').Value = '  -3.18%  '
$ws.Range('D14').Value = '2.760.04'
$ws.Range('E14').Value = '  -19.34%  '
$ws.Range('D15').Value = '54.160.23'
$ws.Range('E15').Value = '  -10.68%  '
$ws.Range('D16').Value = '18.96'
$ws.Range('E16').Value = '  -17.01%  '
$ws.Range('E17').Value = '  -14.27%  '
$ws.Range('D18').Value = '2.369.85'
$ws.Range('E18').Value = '  -18.57%  '
$ws.Range('D19').Value = '4.00'
$ws.Range('E19').Value = '  -20.83%  '
$ws.Range('D20').Value = '299.75'
$ws.Range('E20').Value = '  -17.58%  '
$ws.Range('D21').Value = '9.29'
$ws.Range('E21').Value = '  -20.88%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').Value = '5.62'
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('D24').Value = '5.42'
$ws.Range('E24').Value = '  -18.60%  '
$ws.Range('D25').Value = '55.78'
$ws.Range('E25').Value = '  -14.13%  '
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  -15.21%  '
$ws.Range('D28').Value = '0.370'
$ws.Range('E28').Value = '  -19.01%  '
$ws.Range('D29').Value = '6.98'
$ws.Range('E29').Value = '  -11.67%  '
$ws.Range('D30').Value = '0.996'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').Value = '0.0₃0707'
$ws.Range('E31').Value = '  -18.35%  '
$ws.Range('D32').Value = '147.09'
$ws.Range('E32').Value = '  -3.55%  '
$ws.Range('D33').Value = '17.49'
$ws.Range('E33').Value = '  -11.68%  '
$ws.Range('D34').Value = '1.36'
$ws.Range('E34').Value = '  -19.21%  '
$ws.Range('D35').Value = '4.69'
$ws.Range('E35').Value = '  -16.29%  '
$ws.Range('D36').Value = '3.58'
$ws.Range('E36').Value = '  -18.78%  '
$ws.Range('D37').Value = '0.838'
$ws.Range('E37').Value = '  -17.44%  '
$ws.Range('E38').Value = '  -16.66%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '33.62'
$ws.Range('E39').Value = '  -10.87%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '0.995'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = '10.26'
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').Value = '3.15'
$ws.Range('E42').Value = '  -15.93%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.945.04'
$ws.Range('E43').Value = '  -15.67%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.22'
$ws.Range('E44').Value = '  -18.10%  '
$ws.Range('D45').Value = '0.0497'
$ws.Range('E45').Value = '  -14.89%  '
$ws.Range('D46').Value = '0.531'
$ws.Range('E46').Value = '  -18.33%  '
$ws.Range('D47').Value = '0.0212'
$ws.Range('E47').Value = '  -11.33%  '
$ws.Range('D48').Value = '0.0832'
$ws.Range('E48').Value = '  -10.22%  '
$ws.Range('D49').Value = '16.01'
$ws.Range('E49').Value = '  -22.32%  '
$ws.Range('D50').Value = '4.05'
$ws.Range('E50').Value = '  -19.59%  '
$ws.Range('E51').Value = '  -3.66%  '
